$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 6).Value = 1.76
$ws.Cells.Item(2, 7).Value = 2.24
$ws.Cells.Item(2, 8).Value = 1.8
$ws.Cells.Item(2, 9).Value = 9.800000000000001
$ws.Cells.Item(2, 10).Value = 3.05
$ws.Cells.Item(2, 11).Value = 1000
$ws.Cells.Item(2, 16).Value = 1.74
$ws.Cells.Item(2, 17).Value = 1.8
$ws.Cells.Item(3, 6).Value = 1.04
$ws.Cells.Item(3, 7).Value = 1000
$ws.Cells.Item(3, 8).Value = 1.04
$ws.Cells.Item(3, 9).Value = 1000
$ws.Cells.Item(3, 10).Value = 1.03
$ws.Cells.Item(3, 11).Value = 1000
$ws.Cells.Item(3, 14).Value = 1.25
$ws.Cells.Item(3, 15).Value = 1.01
$ws.Cells.Item(3, 16).Value = 1.25
$ws.Cells.Item(3, 17).Value = 1.34
$ws.Cells.Item(4, 6).Value = 1.32
$ws.Cells.Item(4, 7).Value = 1.52
$ws.Cells.Item(4, 8).Value = 2.92
$ws.Cells.Item(4, 10).Value = 2.92
$ws.Cells.Item(4, 16).Value = 2.38
$ws.Cells.Item(5, 6).Value = 1.26
$ws.Cells.Item(5, 16).Value = 2.56
$ws.Cells.Item(5, 17).Value = 1.44
$ws.Cells.Item(6, 6).Value = 2.42
$ws.Cells.Item(6, 7).Value = 2.5
$ws.Cells.Item(6, 8).Value = 3.45
$ws.Cells.Item(6, 12).Value = 1.01
$ws.Cells.Item(6, 14).Value = 3.2
$ws.Cells.Item(6, 16).Value = 1.77
$ws.Cells.Item(6, 17).Value = 2.24
$ws.Cells.Item(6, 18).Value = 1.28
$ws.Cells.Item(6, 20).Value = 1.91
$ws.Cells.Item(6, 21).Value = 1.95
$ws.Cells.Item(6, 25).Value = 12
$ws.Cells.Item(6, 28).Value = 9.199999999999999
$ws.Cells.Item(6, 31).Value = 50
$ws.Cells.Item(6, 40).Value = 26
$ws.Cells.Item(7, 7).Value = 1.53
$ws.Cells.Item(7, 8).Value = 5
$ws.Cells.Item(7, 9).Value = 9.4
$ws.Cells.Item(7, 10).Value = 4.5
$ws.Cells.Item(7, 16).Value = 2.6
$ws.Cells.Item(8, 6).Value = 1.33
$ws.Cells.Item(8, 7).Value = 1.59
$ws.Cells.Item(8, 8).Value = 2.7
$ws.Cells.Item(8, 9).Value = 1000
$ws.Cells.Item(8, 10).Value = 2.7
$ws.Cells.Item(8, 11).Value = 1000
$ws.Cells.Item(8, 14).Value = 2.08
$ws.Cells.Item(8, 16).Value = 2.08
$ws.Cells.Item(9, 6).Value = 4.8
$ws.Cells.Item(9, 7).Value = 5.5
$ws.Cells.Item(9, 8).Value = 1.81
$ws.Cells.Item(9, 9).Value = 1.95
$ws.Cells.Item(9, 10).Value = 3.55
$ws.Cells.Item(9, 11).Value = 3.95
$ws.Cells.Item(9, 16).Value = 1.78
$ws.Cells.Item(9, 17).Value = 2.02
$ws.Cells.Item(10, 6).Value = 1.04
$ws.Cells.Item(10, 7).Value = 1000
$ws.Cells.Item(10, 8).Value = 1.04
$ws.Cells.Item(10, 9).Value = 1000
$ws.Cells.Item(10, 10).Value = 1.01
$ws.Cells.Item(10, 11).Value = 1000
$ws.Cells.Item(10, 16).Value = 2.48
$ws.Cells.Item(11, 7).Value = 1.56
$ws.Cells.Item(11, 8).Value = 7.4
$ws.Cells.Item(11, 9).Value = 9.800000000000001
$ws.Cells.Item(11, 11).Value = 4.5
$ws.Cells.Item(11, 16).Value = 1.64
$ws.Cells.Item(11, 17).Value = 2.28
$ws.Cells.Item(13, 6).Value = 3.4
$ws.Cells.Item(13, 7).Value = 4.1
$ws.Cells.Item(13, 9).Value = 2.6
$ws.Cells.Item(13, 10).Value = 2.7
$ws.Cells.Item(13, 11).Value = 3.75
$ws.Cells.Item(13, 14).Value = 1.42
$ws.Cells.Item(13, 15).Value = 1.01
$ws.Cells.Item(14, 6).Value = 1.48
$ws.Cells.Item(14, 7).Value = 1.86
$ws.Cells.Item(14, 8).Value = 2.16
$ws.Cells.Item(14, 9).Value = 1000
$ws.Cells.Item(14, 10).Value = 2.16
$ws.Cells.Item(14, 11).Value = 1000
$ws.Cells.Item(14, 14).Value = 1.87
$ws.Cells.Item(14, 15).Value = 1.01
$ws.Cells.Item(14, 17).Value = 1.8
$ws.Cells.Item(15, 6).Value = 1.69
$ws.Cells.Item(15, 8).Value = 1.9
$ws.Cells.Item(15, 9).Value = 2.44
$ws.Cells.Item(15, 10).Value = 1.69
$ws.Cells.Item(15, 17).Value = 2.16
$ws.Cells.Item(16, 6).Value = 1.36
$ws.Cells.Item(16, 7).Value = 1.61
$ws.Cells.Item(16, 8).Value = 2.64
$ws.Cells.Item(16, 10).Value = 2.64
$ws.Cells.Item(16, 11).Value = 1000
$ws.Cells.Item(16, 14).Value = 1.61
$ws.Cells.Item(16, 15).Value = 1.01
$ws.Cells.Item(16, 16).Value = 1.61
$ws.Cells.Item(16, 17).Value = 1.95
$ws.Cells.Item(17, 7).Value = 1.63
$ws.Cells.Item(17, 8).Value = 2.58
$ws.Cells.Item(17, 10).Value = 2.58
$ws.Cells.Item(18, 6).Value = 3.2
$ws.Cells.Item(18, 7).Value = 1000
$ws.Cells.Item(18, 8).Value = 1.4
$ws.Cells.Item(18, 9).Value = 2.46
$ws.Cells.Item(18, 10).Value = 2.94
$ws.Cells.Item(18, 16).Value = 1.58
$ws.Cells.Item(18, 17).Value = 2.04
$ws.Cells.Item(18, 20).Value = 1.01
$ws.Cells.Item(18, 21).Value = 1.01
$ws.Cells.Item(18, 24).Value = 1000
$ws.Cells.Item(18, 25).Value = 1000
$ws.Cells.Item(18, 26).Value = 1000
$ws.Cells.Item(18, 27).Value = 1000
$ws.Cells.Item(18, 28).Value = 1000
$ws.Cells.Item(18, 29).Value = 1000
$ws.Cells.Item(18, 30).Value = 1000
$ws.Cells.Item(18, 31).Value = 1000
$ws.Cells.Item(18, 32).Value = 1000
$ws.Cells.Item(18, 33).Value = 1000
$ws.Cells.Item(18, 34).Value = 1000
$ws.Cells.Item(19, 18).Value = 1.26
$ws.Cells.Item(19, 19).Value = 4.4
$ws.Cells.Item(19, 21).Value = 1.95
$ws.Cells.Item(19, 24).Value = 1000
$ws.Cells.Item(19, 25).Value = 1000
$ws.Cells.Item(19, 26).Value = 1000
$ws.Cells.Item(19, 27).Value = 1000
$ws.Cells.Item(19, 28).Value = 1000
$ws.Cells.Item(19, 30).Value = 1000
$ws.Cells.Item(19, 31).Value = 1000
$ws.Cells.Item(19, 32).Value = 1000
$ws.Cells.Item(19, 33).Value = 12
$ws.Cells.Item(19, 34).Value = 1000
$ws.Cells.Item(19, 35).Value = 1000
$ws.Cells.Item(19, 36).Value = 1000
$ws.Cells.Item(19, 37).Value = 1000
$ws.Cells.Item(19, 38).Value = 1000
$ws.Cells.Item(19, 39).Value = 1000
$ws.Cells.Item(19, 40).Value = 1000
$ws.Cells.Item(20, 6).Value = 1.91
$ws.Cells.Item(20, 7).Value = 1.98
$ws.Cells.Item(20, 8).Value = 4
$ws.Cells.Item(20, 9).Value = 4.9
$ws.Cells.Item(20, 14).Value = 2
$ws.Cells.Item(20, 16).Value = 2
$ws.Cells.Item(20, 17).Value = 1.28
$ws.Cells.Item(20, 18).Value = 1.33
$ws.Cells.Item(20, 19).Value = 2.66
$ws.Cells.Item(20, 20).Value = 1.01
$ws.Cells.Item(20, 21).Value = 1.01
$ws.Cells.Item(20, 24).Value = 1000
$ws.Cells.Item(20, 25).Value = 1000
$ws.Cells.Item(20, 26).Value = 1000
$ws.Cells.Item(20, 27).Value = 1000
$ws.Cells.Item(20, 28).Value = 1000
$ws.Cells.Item(20, 29).Value = 10
$ws.Cells.Item(20, 30).Value = 1000
$ws.Cells.Item(20, 31).Value = 1000
$ws.Cells.Item(20, 32).Value = 1000
$ws.Cells.Item(20, 33).Value = 1000
$ws.Cells.Item(20, 34).Value = 1000
$ws.Cells.Item(20, 35).Value = 1000
$ws.Cells.Item(20, 36).Value = 1000
$ws.Cells.Item(20, 37).Value = 1000
$ws.Cells.Item(20, 38).Value = 1000
$ws.Cells.Item(20, 39).Value = 1000
$ws.Cells.Item(20, 40).Value = 1000
$ws.Cells.Item(20, 41).Value = 1000
$ws.Cells.Item(21, 18).Value = 1.63
$ws.Cells.Item(21, 20).Value = 1.59
$ws.Cells.Item(21, 24).Value = 21
$ws.Cells.Item(21, 26).Value = 44
$ws.Cells.Item(21, 27).Value = 1000
$ws.Cells.Item(21, 29).Value = 9.6
$ws.Cells.Item(21, 31).Value = 50
$ws.Cells.Item(21, 32).Value = 15
$ws.Cells.Item(21, 35).Value = 1000
$ws.Cells.Item(21, 36).Value = 23
$ws.Cells.Item(21, 40).Value = 9
$ws.Cells.Item(21, 41).Value = 50
$ws.Cells.Item(22, 6).Value = 2.66
$ws.Cells.Item(22, 7).Value = 2.74
$ws.Cells.Item(22, 9).Value = 2.9
$ws.Cells.Item(22, 14).Value = 3.8
$ws.Cells.Item(22, 16).Value = 1.92
$ws.Cells.Item(22, 17).Value = 2.04
$ws.Cells.Item(22, 18).Value = 1.35
$ws.Cells.Item(22, 19).Value = 3.7
$ws.Cells.Item(22, 21).Value = 2.18
$ws.Cells.Item(22, 24).Value = 14
$ws.Cells.Item(22, 27).Value = 55
$ws.Cells.Item(22, 31).Value = 34
$ws.Cells.Item(22, 34).Value = 18.5
$ws.Cells.Item(22, 35).Value = 60
$ws.Cells.Item(22, 36).Value = 42
$ws.Cells.Item(22, 37).Value = 32
$ws.Cells.Item(22, 38).Value = 55
$ws.Cells.Item(22, 39).Value = 120
$ws.Cells.Item(22, 41).Value = 30
$ws.Cells.Item(23, 6).Value = 1.19
$ws.Cells.Item(23, 8).Value = 18
$ws.Cells.Item(23, 9).Value = 24
$ws.Cells.Item(23, 11).Value = 8.4
$ws.Cells.Item(23, 16).Value = 2.32
$ws.Cells.Item(23, 17).Value = 1.52
